$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the Price (D) and Volume(1h) (E) columns keep being stored as text,
# matching the source data which writes numeric-looking values (e.g. "1.006",
# "24.230.48") as plain text rather than numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "24.216.32"
$ws.Range("E2").Value = "  +14.49%  "
$ws.Range("D3").Value = "1.671.30"
$ws.Range("E3").Value = "  +8.49%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "307.32"
$ws.Range("E5").Value = "  +8.78%  "
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("D7").Value = "0.3714"
$ws.Range("E7").Value = "  +2.22%  "
$ws.Range("D8").Value = "0.3419"
$ws.Range("E8").Value = "  +7.01%  "
$ws.Range("D9").Value = "47.68"
$ws.Range("E9").Value = "  +16.76%  "
$ws.Range("D10").Value = "1.173"
$ws.Range("E10").Value = "  +6.74%  "
$ws.Range("D11").Value = "0.07253"
$ws.Range("E11").Value = "  +6.26%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "20.36"
$ws.Range("E13").Value = "  +8.65%  "
$ws.Range("D14").Value = "6.062"
$ws.Range("E14").Value = "  +6.55%  "
$ws.Range("D15").Value = "6.729"
$ws.Range("E15").Value = "  +5.55%  "
$ws.Range("D16").Value = "1.676.74"
$ws.Range("E16").Value = "  +9.00%  "
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").Value = "  +5.33%  "
$ws.Range("D18").Value = "0.9994"
$ws.Range("E18").Value = "  +3.29%  "
$ws.Range("D19").Value = "0.06693"
$ws.Range("E19").Value = "  +9.74%  "
$ws.Range("D20").Value = "81.09"
$ws.Range("E20").Value = "  +11.63%  "
$ws.Range("D21").Value = "16.34"
$ws.Range("E21").Value = "  +8.37%  "
$ws.Range("D22").Value = "6.110"
$ws.Range("E22").Value = "  +6.81%  "
$ws.Range("D23").Value = "12.00"
$ws.Range("E23").Value = "  +5.57%  "
$ws.Range("D24").Value = "24.225.92"
$ws.Range("E24").Value = "  +14.34%  "
$ws.Range("D25").Value = "2.397"
$ws.Range("E25").Value = "  +3.18%  "
$ws.Range("D26").Value = "3.361"
$ws.Range("E26").Value = "  -9.28%  "
$ws.Range("D27").Value = "2.633"
$ws.Range("E27").Value = "  +18.26%  "
$ws.Range("D28").Value = "151.70"
$ws.Range("E28").Value = "  +1.97%  "
$ws.Range("D29").Value = "19.41"
$ws.Range("E29").Value = "  +9.75%  "
$ws.Range("D30").Value = "1.862.15"
$ws.Range("E30").Value = "  +9.02%  "
$ws.Range("D31").Value = "126.85"
$ws.Range("E31").Value = "  +7.07%  "
$ws.Range("D32").Value = "6.378"
$ws.Range("E32").Value = "  +22.05%  "
$ws.Range("D33").Value = "4.042"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").Value = "0.9816"
$ws.Range("E34").Value = "  +14.46%  "
$ws.Range("D35").Value = "1.747"
$ws.Range("E35").Value = "  +15.84%  "
$ws.Range("D36").Value = "0.08435"
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("D37").Value = "12.47"
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "0.06381"
$ws.Range("E38").Value = "  +8.87%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.331"
$ws.Range("E39").Value = "  +8.07%  "
$ws.Range("D40").Value = "8.775"
$ws.Range("E40").Value = "  +13.07%  "
$ws.Range("D41").Value = "0.02327"
$ws.Range("E41").Value = "  +10.34%  "
$ws.Range("E42").Value = "  +5.92%  "
$ws.Range("D43").Value = "0.2092"
$ws.Range("E43").Value = "  +8.97%  "
$ws.Range("D44").Value = "0.6122"
$ws.Range("E44").Value = "  +12.20%  "
$ws.Range("D45").Value = "0.9976"
$ws.Range("E45").Value = "  +2.99%  "
$ws.Range("D46").Value = "13.23"
$ws.Range("E46").Value = "  +5.06%  "
$ws.Range("D47").Value = "3.796"
$ws.Range("E47").Value = "  +6.27%  "
$ws.Range("D48").Value = "0.5918"
$ws.Range("E48").Value = "  +8.73%  "
$ws.Range("D49").Value = "127.60"
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("D50").Value = "2.006"
$ws.Range("E50").Value = "  +7.06%  "
$ws.Range("D51").Value = "0.07154"
$ws.Range("E51").Value = "  +8.69%  "
